# Generate Report for Handback
#
# This script updates the localization-status workbook so that the two
# files "53169d0f-da5d-4779-bcd7-24a9315d6594.md" and
# "75e7169d-b8e4-4486-8015-6844cff8a982.md" are reported as handed back
# (in sync with en-US), populates their "Latest Target File" / "Latest
# Handback File" / "Latest Handback DateTime" columns on the per-locale
# sheets, and re-sorts every sheet by file name (53169d0f, 75e7169d,
# eba35525, e34a3baa).

$wb = $excel.ActiveWorkbook

$HANDED_BACK = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Clear every hyperlink on the sheet first (Range.Hyperlinks.Delete()
# removes all hyperlinks that belong to the worksheet in this host) so we
# can rebuild them cleanly in the new row order.
$ov.Range("A1").Hyperlinks.Delete()

$ov.Range("A2").Value = "53169d0f-da5d-4779-bcd7-24a9315d6594.md"
$ov.Range("B2").Value = $HANDED_BACK
$ov.Range("C2").Value = $HANDED_BACK
$ov.Range("D2").Value = "2016-03-24 00:19:42"

$ov.Range("A3").Value = "75e7169d-b8e4-4486-8015-6844cff8a982.md"
$ov.Range("B3").Value = $HANDED_BACK
$ov.Range("C3").Value = $HANDED_BACK
$ov.Range("D3").Value = "2016-03-24 00:19:42"

$ov.Range("A4").Value = "eba35525-9799-416b-a933-f6301e488292.md"
$ov.Range("B4").Value = "In Translation"
$ov.Range("C4").Value = "In Translation"
$ov.Range("D4").Value = "2016-03-24 00:18:11"

$ov.Range("A5").Value = "e34a3baa-cbab-416f-b3a8-a96cf004f085.md"
$ov.Range("B5").Value = "Ready for handoff"
$ov.Range("C5").Value = "Ready for handoff"
$ov.Range("D5").Value = "2016-03-24 00:19:42"

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d0b8de1f10223117fb73313fcba12cf2553f756a/e2e/53169d0f-da5d-4779-bcd7-24a9315d6594.md", "", "", "53169d0f-da5d-4779-bcd7-24a9315d6594.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d0b8de1f10223117fb73313fcba12cf2553f756a/e2e/75e7169d-b8e4-4486-8015-6844cff8a982.md", "", "", "75e7169d-b8e4-4486-8015-6844cff8a982.md")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f751a3118795b42ada48d1d4e01aa429e08a8efb/e2e/eba35525-9799-416b-a933-f6301e488292.md", "", "", "eba35525-9799-416b-a933-f6301e488292.md")
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/beb454a33aef493a34ef277aedb840958e7ed78e/e2e/e34a3baa-cbab-416f-b3a8-a96cf004f085.md", "", "", "e34a3baa-cbab-416f-b3a8-a96cf004f085.md")

# Re-apply the values (Hyperlinks.Add replaces the cell text with the
# TextToDisplay argument, which is the same string here, but keep this
# for safety/robustness across hosts).
$ov.Range("A2").Value = "53169d0f-da5d-4779-bcd7-24a9315d6594.md"
$ov.Range("A3").Value = "75e7169d-b8e4-4486-8015-6844cff8a982.md"
$ov.Range("A4").Value = "eba35525-9799-416b-a933-f6301e488292.md"
$ov.Range("A5").Value = "e34a3baa-cbab-416f-b3a8-a96cf004f085.md"

# ---------------------------------------------------------------------
# Helper data for the per-locale sheets (zh-cn / de-de)
# ---------------------------------------------------------------------
# Column layout (A..L):
#  A Source File Name        G Latest Handback File
#  B File Extension          H Latest Handback DateTime
#  C Status                  I Reference Tokens
#  D Latest Handoff File     J Handoff Reason
#  E Latest Handoff Datetime K Dependency From
#  F Latest Target File      L Error Detail

function Set-LocaleSheet($sheet, $locale, $row2HandoffDate, $row3HandoffDate, $row4HandoffDate, $row5HandoffDate, $handbackDateTime) {

    # Clear all existing hyperlinks on this sheet before rebuilding them.
    $sheet.Range("A1").Hyperlinks.Delete()

    # --- Row 2: 53169d0f... (handed back) ---
    $sheet.Range("A2").Value = "53169d0f-da5d-4779-bcd7-24a9315d6594.md"
    $sheet.Range("B2").Value = ".md"
    $sheet.Range("C2").Value = $HANDED_BACK
    $sheet.Range("D2").Value = "53169d0f-da5d-4779-bcd7-24a9315d6594.661293e2b235a2848eac94d2a084b67bb016dc17.$locale.xlf"
    $sheet.Range("E2").Value = $row2HandoffDate
    $sheet.Range("F2").Value = "53169d0f-da5d-4779-bcd7-24a9315d6594.md"
    $sheet.Range("G2").Value = "53169d0f-da5d-4779-bcd7-24a9315d6594.661293e2b235a2848eac94d2a084b67bb016dc17.$locale.xlf"
    $sheet.Range("H2").Value = $handbackDateTime
    $sheet.Range("J2").Value = "Include"

    # --- Row 3: 75e7169d... (handed back) ---
    $sheet.Range("A3").Value = "75e7169d-b8e4-4486-8015-6844cff8a982.md"
    $sheet.Range("B3").Value = ".md"
    $sheet.Range("C3").Value = $HANDED_BACK
    $sheet.Range("D3").Value = "75e7169d-b8e4-4486-8015-6844cff8a982.d0daae3d6ad2f49b2c2e10336dc9ac4b72f757b1.$locale.xlf"
    $sheet.Range("E3").Value = $row3HandoffDate
    $sheet.Range("F3").Value = "75e7169d-b8e4-4486-8015-6844cff8a982.md"
    $sheet.Range("G3").Value = "75e7169d-b8e4-4486-8015-6844cff8a982.d0daae3d6ad2f49b2c2e10336dc9ac4b72f757b1.$locale.xlf"
    $sheet.Range("H3").Value = $handbackDateTime
    $sheet.Range("J3").Value = "Include"

    # --- Row 4: eba35525... (in translation) ---
    $sheet.Range("A4").Value = "eba35525-9799-416b-a933-f6301e488292.md"
    $sheet.Range("B4").Value = ".md"
    $sheet.Range("C4").Value = "In Translation"
    $sheet.Range("D4").Value = "eba35525-9799-416b-a933-f6301e488292.2841483b9fc535cc2216d4f730eb1c5fe309e396.$locale.xlf"
    $sheet.Range("E4").Value = $row4HandoffDate
    $sheet.Range("F4").Value = ""
    $sheet.Range("G4").Value = ""
    $sheet.Range("H4").Value = "0001-01-01 00:00:00"
    $sheet.Range("J4").Value = "Include"

    # --- Row 5: e34a3baa... (ready for handoff) ---
    $sheet.Range("A5").Value = "e34a3baa-cbab-416f-b3a8-a96cf004f085.md"
    $sheet.Range("B5").Value = ".md"
    $sheet.Range("C5").Value = "Ready for handoff"
    $sheet.Range("D5").Value = "e34a3baa-cbab-416f-b3a8-a96cf004f085.cba5542ff95e188e91e448ff331d5d1df50a40da.$locale.xlf"
    $sheet.Range("E5").Value = $row5HandoffDate
    $sheet.Range("F5").Value = ""
    $sheet.Range("G5").Value = ""
    $sheet.Range("H5").Value = "0001-01-01 00:00:00"
    $sheet.Range("J5").Value = "Include"

    # --- Hyperlinks ---
    $sheet.Hyperlinks.Add($sheet.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d0b8de1f10223117fb73313fcba12cf2553f756a/e2e/53169d0f-da5d-4779-bcd7-24a9315d6594.md", "", "", "53169d0f-da5d-4779-bcd7-24a9315d6594.md")
    $sheet.Hyperlinks.Add($sheet.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8961b0e3df67928226b5d419f16c6d290f5b4da/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/53169d0f-da5d-4779-bcd7-24a9315d6594.661293e2b235a2848eac94d2a084b67bb016dc17.$locale.xlf", "", "", "53169d0f-da5d-4779-bcd7-24a9315d6594.661293e2b235a2848eac94d2a084b67bb016dc17.$locale.xlf")
    $sheet.Hyperlinks.Add($sheet.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/d0b8de1f10223117fb73313fcba12cf2553f756a/e2e/53169d0f-da5d-4779-bcd7-24a9315d6594.md", "", "", "53169d0f-da5d-4779-bcd7-24a9315d6594.md")
    $sheet.Hyperlinks.Add($sheet.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8961b0e3df67928226b5d419f16c6d290f5b4da/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/53169d0f-da5d-4779-bcd7-24a9315d6594.661293e2b235a2848eac94d2a084b67bb016dc17.$locale.xlf", "", "", "53169d0f-da5d-4779-bcd7-24a9315d6594.661293e2b235a2848eac94d2a084b67bb016dc17.$locale.xlf")

    $sheet.Hyperlinks.Add($sheet.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d0b8de1f10223117fb73313fcba12cf2553f756a/e2e/75e7169d-b8e4-4486-8015-6844cff8a982.md", "", "", "75e7169d-b8e4-4486-8015-6844cff8a982.md")
    $sheet.Hyperlinks.Add($sheet.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8961b0e3df67928226b5d419f16c6d290f5b4da/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/75e7169d-b8e4-4486-8015-6844cff8a982.d0daae3d6ad2f49b2c2e10336dc9ac4b72f757b1.$locale.xlf", "", "", "75e7169d-b8e4-4486-8015-6844cff8a982.d0daae3d6ad2f49b2c2e10336dc9ac4b72f757b1.$locale.xlf")
    $sheet.Hyperlinks.Add($sheet.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/d0b8de1f10223117fb73313fcba12cf2553f756a/e2e/75e7169d-b8e4-4486-8015-6844cff8a982.md", "", "", "75e7169d-b8e4-4486-8015-6844cff8a982.md")
    $sheet.Hyperlinks.Add($sheet.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8961b0e3df67928226b5d419f16c6d290f5b4da/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/75e7169d-b8e4-4486-8015-6844cff8a982.d0daae3d6ad2f49b2c2e10336dc9ac4b72f757b1.$locale.xlf", "", "", "75e7169d-b8e4-4486-8015-6844cff8a982.d0daae3d6ad2f49b2c2e10336dc9ac4b72f757b1.$locale.xlf")

    $sheet.Hyperlinks.Add($sheet.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f751a3118795b42ada48d1d4e01aa429e08a8efb/e2e/eba35525-9799-416b-a933-f6301e488292.md", "", "", "eba35525-9799-416b-a933-f6301e488292.md")
    $sheet.Hyperlinks.Add($sheet.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9716e3fe034852508e4df56c9615a2048659a4a5/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht/eba35525-9799-416b-a933-f6301e488292.2841483b9fc535cc2216d4f730eb1c5fe309e396.$locale.xlf", "", "", "eba35525-9799-416b-a933-f6301e488292.2841483b9fc535cc2216d4f730eb1c5fe309e396.$locale.xlf")

    $sheet.Hyperlinks.Add($sheet.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/beb454a33aef493a34ef277aedb840958e7ed78e/e2e/e34a3baa-cbab-416f-b3a8-a96cf004f085.md", "", "", "e34a3baa-cbab-416f-b3a8-a96cf004f085.md")
    $sheet.Hyperlinks.Add($sheet.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c8961b0e3df67928226b5d419f16c6d290f5b4da/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/e34a3baa-cbab-416f-b3a8-a96cf004f085.cba5542ff95e188e91e448ff331d5d1df50a40da.$locale.xlf", "", "", "e34a3baa-cbab-416f-b3a8-a96cf004f085.cba5542ff95e188e91e448ff331d5d1df50a40da.$locale.xlf")

    # Re-apply text values, since Hyperlinks.Add overwrites the cell text
    # with TextToDisplay (identical values here, kept for robustness).
    $sheet.Range("A2").Value = "53169d0f-da5d-4779-bcd7-24a9315d6594.md"
    $sheet.Range("D2").Value = "53169d0f-da5d-4779-bcd7-24a9315d6594.661293e2b235a2848eac94d2a084b67bb016dc17.$locale.xlf"
    $sheet.Range("F2").Value = "53169d0f-da5d-4779-bcd7-24a9315d6594.md"
    $sheet.Range("G2").Value = "53169d0f-da5d-4779-bcd7-24a9315d6594.661293e2b235a2848eac94d2a084b67bb016dc17.$locale.xlf"

    $sheet.Range("A3").Value = "75e7169d-b8e4-4486-8015-6844cff8a982.md"
    $sheet.Range("D3").Value = "75e7169d-b8e4-4486-8015-6844cff8a982.d0daae3d6ad2f49b2c2e10336dc9ac4b72f757b1.$locale.xlf"
    $sheet.Range("F3").Value = "75e7169d-b8e4-4486-8015-6844cff8a982.md"
    $sheet.Range("G3").Value = "75e7169d-b8e4-4486-8015-6844cff8a982.d0daae3d6ad2f49b2c2e10336dc9ac4b72f757b1.$locale.xlf"

    $sheet.Range("A4").Value = "eba35525-9799-416b-a933-f6301e488292.md"
    $sheet.Range("D4").Value = "eba35525-9799-416b-a933-f6301e488292.2841483b9fc535cc2216d4f730eb1c5fe309e396.$locale.xlf"

    $sheet.Range("A5").Value = "e34a3baa-cbab-416f-b3a8-a96cf004f085.md"
    $sheet.Range("D5").Value = "e34a3baa-cbab-416f-b3a8-a96cf004f085.cba5542ff95e188e91e448ff331d5d1df50a40da.$locale.xlf"
}

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
Set-LocaleSheet $zh "zh-cn" `
    "2016-03-24 00:19:37" `
    "2016-03-24 00:19:37" `
    "2016-03-24 00:18:07" `
    "2016-03-24 00:19:37" `
    "2016-03-24 00:20:01"

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
Set-LocaleSheet $de "de-de" `
    "2016-03-24 00:19:42" `
    "2016-03-24 00:19:42" `
    "2016-03-24 00:18:11" `
    "2016-03-24 00:19:42" `
    "2016-03-24 00:20:10"
